$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.710.93"
$ws.Range("E2").Value = "  +10.16%  "
# Row 3
$ws.Range("D3").Value = "2.683.97"
$ws.Range("E3").Value = "  +14.63%  "
# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.17%  "
# Row 5
$ws.Range("D5").Value = "'515.50"
$ws.Range("E5").Value = "  +9.02%  "
# Row 6
$ws.Range("D6").Value = "'161.69"
$ws.Range("E6").Value = "  +12.37%  "
# Row 7
$ws.Range("D7").Value = "'0.614"
$ws.Range("E7").Value = "  +5.27%  "
# Row 8
$ws.Range("D8").Value = "'0.996"
$ws.Range("E8").Value = "  -0.18%  "
# Row 9
$ws.Range("D9").Value = "2.687.02"
$ws.Range("E9").Value = "  +14.65%  "
# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.107"
$ws.Range("E10").Value = "  +13.34%  "
# Row 11
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "'6.14"
$ws.Range("E11").Value = "  +13.01%  "
# Row 12
$ws.Range("D12").Value = "'0.352"
$ws.Range("E12").Value = "  +9.15%  "
# Row 13
$ws.Range("E13").Value = "  +2.11%  "
# Row 14
$ws.Range("D14").Value = "3.144.12"
$ws.Range("E14").Value = "  +14.86%  "
# Row 15
$ws.Range("D15").Value = "61.021.25"
$ws.Range("E15").Value = "  +10.65%  "
# Row 16
$ws.Range("D16").Value = "'22.67"
$ws.Range("E16").Value = "  +14.69%  "
# Row 17
$ws.Range("D17").Value = "'0.0000142"
$ws.Range("E17").Value = "  +11.00%  "
# Row 18
$ws.Range("D18").Value = "2.684.36"
$ws.Range("E18").Value = "  +14.86%  "
# Row 19
$ws.Range("D19").Value = "'4.87"
$ws.Range("E19").Value = "  +7.68%  "
# Row 20
$ws.Range("D20").Value = "'356.17"
$ws.Range("E20").Value = "  +13.91%  "
# Row 21
$ws.Range("D21").Value = "'10.66"
$ws.Range("E21").Value = "  +12.63%  "
# Row 22
$ws.Range("D22").Value = "'6.25"
$ws.Range("E22").Value = "  +10.52%  "
# Row 23
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.06%  "
# Row 24
$ws.Range("D24").Value = "'60.65"
$ws.Range("E24").Value = "  +9.13%  "
# Row 25
$ws.Range("D25").Value = "'0.430"
$ws.Range("E25").Value = "  +9.96%  "
# Row 26
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.790.91"
$ws.Range("E26").Value = "  +14.66%  "
# Row 27
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.169"
$ws.Range("E27").Value = "  +10.96%  "
# Row 28
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.03%  "
# Row 29
$ws.Range("D29").Value = "0.0₃0892"
$ws.Range("E29").Value = "  +22.04%  "
# Row 30
$ws.Range("D30").Value = "'7.65"
$ws.Range("E30").Value = "  +9.06%  "
# Row 31
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.07%  "
# Row 32
$ws.Range("D32").Value = "'19.85"
$ws.Range("E32").Value = "  +9.78%  "
# Row 33
$ws.Range("D33").Value = "'158.75"
$ws.Range("E33").Value = "  +8.72%  "
# Row 34
$ws.Range("D34").Value = "'1.60"
$ws.Range("E34").Value = "  +9.11%  "
# Row 35
$ws.Range("D35").Value = "'5.75"
$ws.Range("E35").Value = "  +13.79%  "
# Row 36
$ws.Range("D36").Value = "'4.10"
$ws.Range("E36").Value = "  +16.73%  "
# Row 37
$ws.Range("D37").Value = "'1.24"
$ws.Range("E37").Value = "  +14.51%  "
# Row 38
$ws.Range("E38").Value = "  +11.24%  "
# Row 39
$ws.Range("D39").Value = "'1.55"
$ws.Range("E39").Value = "  +18.08%  "
# Row 40
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'3.84"
$ws.Range("E40").Value = "  +14.99%  "
# Row 41
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").Value = "'0.849"
$ws.Range("E41").Value = "  +38.84%  "
# Row 42
$ws.Range("D42").Value = "'301.93"
$ws.Range("E42").Value = "  +22.00%  "
# Row 43
$ws.Range("D43").Value = "'36.09"
$ws.Range("E43").Value = "  +7.78%  "
# Row 44
$ws.Range("D44").Value = "'0.650"
$ws.Range("E44").Value = "  +13.72%  "
# Row 45
$ws.Range("D45").Value = "'0.0583"
$ws.Range("E45").Value = "  +13.81%  "
# Row 46
$ws.Range("D46").Value = "'0.102"
$ws.Range("E46").Value = "  +3.42%  "
# Row 47
$ws.Range("D47").Value = "'20.31"
$ws.Range("E47").Value = "  +23.04%  "
# Row 48
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'0.994"
$ws.Range("E48").Value = "  -0.35%  "
# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'5.05"
$ws.Range("E49").Value = "  +17.35%  "
# Row 50
$ws.Range("D50").Value = "'0.0241"
$ws.Range("E50").Value = "  +10.32%  "
# Row 51
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.059.04"
$ws.Range("E51").Value = "  +14.67%  "
